# =====================================================================
# Adds two new laptop entries ("TUF F16 FX608LP" and "XGM 蛟龍16Pro")
# to the SPEC / CPU / GPU / FPS_QHD sheets (rows 20-21), tweaks two
# existing SPEC values (AI15, U19), and restores the cursor / selection
# left behind in each sheet after the edit.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# SPEC sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SPEC")

# two isolated value tweaks on existing rows
$ws.Range("AI15").Value = 2.44
$ws.Range("U19").Value = 2

# new row 20 - TUF F16 FX608LP
$ws.Range("A20").Value = "TUF F16 FX608LP"
$ws.Range("B20").Value = "Intel"
$ws.Range("C20").Value = "U7 255HX"
$ws.Range("D20").Value = 135
$ws.Range("E20").Value = 90
$ws.Range("F20").Value = "RTX5070"
$ws.Range("G20").Value = "GDDR7 8GB"
$ws.Range("H20").Value = "115W"
$ws.Range("I20").Value = 50
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 170
$ws.Range("L20").Value = 16
$ws.Range("M20").Value = "IPS"
$ws.Range("N20").Value = "2560x1600"
$ws.Range("O20").Value = 165
$ws.Range("P20").Value = 400
$ws.Range("Q20").Value = "Y"
$ws.Range("R20").Value = "DDR5 5600MHz"
$ws.Range("S20").Value = "-"
$ws.Range("T20").Value = 2
$ws.Range("U20").Value = 2
$ws.Range("V20").Value = "2*G4x4"
$ws.Range("W20").Value = 1
$ws.Range("X20").Value = "WIFI 6E/BT 5.3"
$ws.Range("Y20").Value = "2A2C"
$ws.Range("Z20").Value = "-"
$ws.Range("AA20").Value = 1
$ws.Range("AB20").Value = 100
$ws.Range("AC20").Value = 2
$ws.Range("AD20").Value = "FHD IR"
$ws.Range("AE20").Value = 90
$ws.Range("AF20").Value = 280
$ws.Range("AG20").Value = "354x269"
$ws.Range("AH20").Value = "17.9~27.3"
$ws.Range("AI20").Value = 2.2

# new row 21 - XGM 蛟龍16Pro
$ws.Range("A21").Value = "XGM 蛟龍16Pro"
$ws.Range("B21").Value = "AMD"
$ws.Range("C21").Value = "R9 9955HX"
$ws.Range("D21").Value = 110
$ws.Range("E21").Value = 110
$ws.Range("F21").Value = "RTX5070Ti"
$ws.Range("G21").Value = "GDDR7 12GB"
$ws.Range("H21").Value = "140W"
$ws.Range("I21").Value = "-"
$ws.Range("J21").Value = "-"
$ws.Range("K21").Value = 205
$ws.Range("L21").Value = 16
$ws.Range("M21").Value = "IPS"
$ws.Range("N21").Value = "2560x1600"
$ws.Range("O21").Value = 300
$ws.Range("P21").Value = 500
$ws.Range("Q21").Value = "Y"
$ws.Range("R21").Value = "DDR5 5600MHz"
$ws.Range("S21").Value = "-"
$ws.Range("T21").Value = 2
$ws.Range("U21").Value = 2
$ws.Range("V21").Value = "2*G4x4"
$ws.Range("W21").Value = 1
$ws.Range("X21").Value = "WIFI 6E/BT 5.3"
$ws.Range("Y21").Value = "3A2C"
$ws.Range("Z21").Value = "-"
$ws.Range("AA21").Value = "-"
$ws.Range("AB21").Value = 100
$ws.Range("AC21").Value = 2
$ws.Range("AD21").Value = "720P IR"
$ws.Range("AE21").Value = 80
$ws.Range("AF21").Value = 280
$ws.Range("AG21").Value = "356.68×253.8"
$ws.Range("AH21").Value = 24.8
$ws.Range("AI21").Value = 2.44

$ws.Range("C30").Select()

# ---------------------------------------------------------------------
# CPU sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CPU")

$ws.Range("A20").Value = "TUF F16 FX608LP"
$ws.Range("B20").Value = "U7 255HX"
$ws.Range("C20").Value = "RTX5070"
$ws.Range("D20").Value = "115W"
$ws.Range("E20").Value = 1959
$ws.Range("G20").Value = 27843

$ws.Range("A21").Value = "XGM 蛟龍16Pro"
$ws.Range("B21").Value = "R9 9955HX"
$ws.Range("C21").Value = "RTX5070Ti"
$ws.Range("D21").Value = "140W"
$ws.Range("E21").Value = 2126
$ws.Range("G21").Value = 37992

# Row 21's C/D cells use the plain "Normal" look (like rows 18/19)
# instead of the column's default bold style - copy that formatting
# over from row 18 without touching its content.
$ws.Range("C18:D18").Copy()
$ws.Range("C21:D21").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(21).RowHeight = 15.5

$ws.Range("E21").Select()

# ---------------------------------------------------------------------
# GPU sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GPU")

$ws.Range("A20").Value = "TUF F16 FX608LP"
$ws.Range("B20").Value = "U7 255HX"
$ws.Range("C20").Value = "RTX5070"
$ws.Range("D20").Value = "115W"
$ws.Range("H20").Value = 14120

$ws.Range("A21").Value = "XGM 蛟龍16Pro"
$ws.Range("B21").Value = "R9 9955HX"
$ws.Range("C21").Value = "RTX5070Ti"
$ws.Range("D21").Value = "140W"
$ws.Range("H21").Value = 17652

# Row 21's C/D cells use the plain "Normal" look (like row 18) instead
# of the column's default bold style.
$ws.Range("C18:D18").Copy()
$ws.Range("C21:D21").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("H22").Select()

# ---------------------------------------------------------------------
# FPS_QHD sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FPS_QHD")

$ws.Range("A20").Value = "TUF F16 FX608LP"
$ws.Range("B20").Value = "U7 255HX"
$ws.Range("C20").Value = "RTX5070"
$ws.Range("D20").Value = "115W"
$ws.Range("F20").Value = 76
$ws.Range("K20").Value = 123

$ws.Range("A21").Value = "XGM 蛟龍16Pro"
$ws.Range("B21").Value = "R9 9955HX"
$ws.Range("C21").Value = "RTX5070Ti"
$ws.Range("D21").Value = "140W"
$ws.Range("F21").Value = 129
$ws.Range("K21").Value = 156

# Row 20's C/D cells use the bold "MGP" style (like rows 2-17) instead
# of the column's plain default style.
$ws.Range("C17:D17").Copy()
$ws.Range("C20:D20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("E30").Select()

# ---------------------------------------------------------------------
# FPS_FHD sheet - selection only, no data changes
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FPS_FHD")
$ws.Range("D28").Select()

# Leave SPEC as the active sheet/tab, matching the saved file.
$wb.Worksheets.Item("SPEC").Activate()
